# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Net data change: the "Periodo Mora" (column E) and "Valor Mora" (column F)
# values for rows 16 and 17 are swapped:
#   Row16: Periodo 2102 / Valor 16959  ->  Periodo 2103 / Valor 36341
#   Row17: Periodo 2103 / Valor 36341  ->  Periodo 2102 / Valor 16959

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2103"
$ws.Range("F16").Value = 36341

$ws.Range("E17").Value = "2102"
$ws.Range("F17").Value = 16959
